$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# The "FilesTab" row's query (B4) is replaced with a revised Cypher query
# that drops the "File Type" and "Breed" output columns.
$newFileQuery = @'

MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
 MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
WHERE diag.stage_of_disease IN ['IVa']
WITH DISTINCT f, parent, c, demo, diag, s
RETURN coalesce(f.file_name, '') AS `File Name`, 
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

$ws.Range("B4").Value = $newFileQuery

# Update the saved selection / active cell to B4 (matches the author's
# last on-screen selection when the file was saved).
$ws.Range("B4").Select() | Out-Null
